# [BI-1158] update test and remove abbreviations from trait table
#
# The "Template" sheet's trait table drops the "Trait abbreviations" column
# (column B) entirely, shifting every later column one to the left, and a
# handful of the remaining header cells are renamed to shorter labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Activate()

# Remove the "Trait abbreviations" column (old column B) - this shifts
# everything from old column C onward one column to the left.
$ws.Columns("B").Delete()

# Rename header cells (columns already shifted left by the delete above):
#   A1 "Ontology term name" -> "Name"
#   B1 "Trait synonyms"     -> "Synonyms"   (was old column C)
#   C1 "Trait description"  -> "Description" (was old column D)
#   F1 "Trait status"       -> "Status"     (was old column G)
#   L1 "Scale name"         -> "Units"      (was old column M)
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Synonyms"
$ws.Range("C1").Value = "Description"
$ws.Range("F1").Value = "Status"
$ws.Range("L1").Value = "Units"

# Update the active selection to reflect where the editor left off.
$ws.Range("L1").Select()
